$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "Cálculo de una variable- Trascendentes y tempranas `$ STEWART, James"
$ws.Range("C6").Value = "Precálculo - Matemáticas para el cálculo `$ STEWART, James"
$ws.Range("C7").Value = "Cálculo de una variable- Trascendentes y tempranas `$ STEWART, James"
$ws.Range("C8").Value = "Precálculo - Matemáticas para el cálculo `$ STEWART, James"
$ws.Range("C9").Value = "Cálculo de una variable- Trascendentes y tempranas `$ STEWART, James"
$ws.Range("C10").Value = "Precálculo - Matemáticas para el cálculo `$ STEWART, James"

$ws.Range("C5:C10").Select()
